$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

$ws.Range("A1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)

# Updated metric values on row 2
$ws.Range("B2").Value = 0.1270053685299624
$ws.Range("C2").Value = 0.9906163806958945
$ws.Range("D2").Value = 0.2768938502609059

# Updated model description text
$ws.Range("F2").Value = "Pipeline(steps=[('model', RandomForestRegressor(max_depth=5))])"

# New data cells
$ws.Range("G2").Value = 0.1434780816666413
$ws.Range("H2").Value = 0.992
